# Auto-generated edit script: applies the cell-value changes described
# by the Jenova_Profits.xlsx diff (profit-tracking workbook, 8 job sheets).
# Every touched cell is a plain numeric literal (no formulas in this file),
# so we set .Value directly; two cells are removed entirely in the target
# (diff drops the <c> element) and are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 25901.076
$ws.Range("I8").Value = 111246.78
$ws.Range("K8").Value = 333740.34
$ws.Range("M8").Value = -333601.34
$ws.Range("H12").Value = 5335.8335
$ws.Range("J12").Value = 7190
$ws.Range("L12").Value = 7190
$ws.Range("N12").Value = -7530
$ws.Range("H13").Value = 5333.3335
$ws.Range("J13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3338
$ws.Range("H17").Value = 33928.395
$ws.Range("J17").Value = 34473.867
$ws.Range("L17").Value = 103421.601
$ws.Range("N17").Value = -103757.601
$ws.Range("H31").Value = 376.4
$ws.Range("I31").Value = 220.5
$ws.Range("K31").Value = 661.5
$ws.Range("M31").Value = -431.5
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1674
$ws.Range("H33").Value = 826.1818
$ws.Range("I33").Value = 888.6
$ws.Range("K33").Value = 888.6
$ws.Range("M33").Value = -659.6
$ws.Range("H43").Value = 3679.8667
$ws.Range("I43").Value = 3781.25
$ws.Range("J43").Value = 3564
$ws.Range("K43").Value = 3781.25
$ws.Range("L43").Value = 3564
$ws.Range("M43").Value = -3712.25
$ws.Range("N43").Value = -3702
$ws.Range("H69").Value = 7992.25
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25748
$ws.Range("H72").Value = 7992.25
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80736
$ws.Range("H80").Value = 2334
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 7500
$ws.Range("N80").Value = -9496
$ws.Range("H83").Value = 2334
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 22500
$ws.Range("N83").Value = -32484
$ws.Range("H92").Value = 538.06665
$ws.Range("I92").Value = 39.272728
$ws.Range("K92").Value = 39.272728
$ws.Range("M92").Value = 1208.727272
$ws.Range("H98").Value = 1959.2667
$ws.Range("I98").Value = 990.75
$ws.Range("K98").Value = 990.75
$ws.Range("M98").Value = 507.25
$ws.Range("H106").Value = 3302.2307
$ws.Range("I106").Value = 5659.75
$ws.Range("K106").Value = 5659.75
$ws.Range("M106").Value = -5028.75
$ws.Range("H107").Value = 36309.715
$ws.Range("I107").Value = 50565.8
$ws.Range("J107").Value = 669.5
$ws.Range("K107").Value = 50565.8
$ws.Range("L107").Value = 669.5
$ws.Range("M107").Value = -48645.8
$ws.Range("N107").Value = -4509.5
$ws.Range("H111").Value = 68545.336
$ws.Range("I111").Value = 92382.55
$ws.Range("K111").Value = 277147.65
$ws.Range("M111").Value = -274080.65
$ws.Range("H112").Value = 3015.1667
$ws.Range("J112").Value = 3072.9312
$ws.Range("L112").Value = 9218.793600000001
$ws.Range("N112").Value = -11434.7936
$ws.Range("H113").Value = 4999.5
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -1745
$ws.Range("H116").Value = 3997.75
$ws.Range("I116").Value = 3856.8572
$ws.Range("J116").Value = 4195
$ws.Range("K116").Value = 3856.8572
$ws.Range("L116").Value = 4195
$ws.Range("M116").Value = -414.8571999999999
$ws.Range("N116").Value = -11079
$ws.Range("H122").Value = 1959.2667
$ws.Range("I122").Value = 990.75
$ws.Range("K122").Value = 2972.25
$ws.Range("M122").Value = -522.25
$ws.Range("H124").Value = 75416.664
$ws.Range("J124").Value = 75416.664
$ws.Range("L124").Value = 75416.664
$ws.Range("N124").Value = -85236.664
$ws.Range("H127").Value = 21867
$ws.Range("I127").Value = 23705.143
$ws.Range("J127").Value = 9000
$ws.Range("K127").Value = 71115.429
$ws.Range("L127").Value = 27000
$ws.Range("M127").Value = -66155.429
$ws.Range("N127").Value = -36920
$ws.Range("H129").Value = 16614.4
$ws.Range("I129").Value = 28717.75
$ws.Range("K129").Value = 86153.25
$ws.Range("M129").Value = -81153.25
$ws.Range("H132").Value = 2490.634
$ws.Range("I132").Value = 1427.7576
$ws.Range("K132").Value = 4283.2728
$ws.Range("M132").Value = -1753.2728
$ws.Range("H133").Value = 58696.938
$ws.Range("J133").Value = 58696.938
$ws.Range("L133").Value = 58696.938
$ws.Range("N133").Value = -68816.93799999999
$ws.Range("H137").Value = 3376.3333
$ws.Range("I137").Value = 1747.421
$ws.Range("J137").Value = 4721.9565
$ws.Range("K137").Value = 5242.263
$ws.Range("L137").Value = 14165.8695
$ws.Range("M137").Value = -2692.263
$ws.Range("N137").Value = -19265.8695
$ws.Range("H138").Value = 4588.7285
$ws.Range("J138").Value = 6447.561
$ws.Range("L138").Value = 19342.683
$ws.Range("N138").Value = -29622.683

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 2166.8
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 2633.5
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 2633.5
$ws.Range("M17").Value = -127
$ws.Range("N17").Value = -2979.5
$ws.Range("H32").Value = 3825.7358
$ws.Range("I32").Value = 3321.8367
$ws.Range("K32").Value = 3321.8367
$ws.Range("M32").Value = -3034.8367
$ws.Range("H45").Value = 10725.75
$ws.Range("I45").Value = 4826.5
$ws.Range("J45").Value = 16625
$ws.Range("K45").Value = 4826.5
$ws.Range("L45").Value = 16625
$ws.Range("M45").Value = -4449.5
$ws.Range("N45").Value = -17379
$ws.Range("H63").Value = 3149.1
$ws.Range("J63").Value = 3401
$ws.Range("L63").Value = 3401
$ws.Range("N63").Value = -4773
$ws.Range("H66").Value = 3149.1
$ws.Range("J66").Value = 3401
$ws.Range("L66").Value = 17005
$ws.Range("N66").Value = -23869
$ws.Range("H92").Value = 75000000
$ws.Range("J92").Value = 100000000
$ws.Range("L92").Value = 100000000
$ws.Range("N92").Value = -100004992
$ws.Range("H102").Value = 2290.6296
$ws.Range("J102").Value = 1031.2
$ws.Range("L102").Value = 1031.2
$ws.Range("N102").Value = -4275.2
$ws.Range("H122").Value = 4199.8667
$ws.Range("I122").Value = 3405
$ws.Range("K122").Value = 10215
$ws.Range("M122").Value = -7765
$ws.Range("H132").Value = 4082.6406
$ws.Range("I132").Value = 1275.7441
$ws.Range("K132").Value = 3827.2323
$ws.Range("M132").Value = -1297.2323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1735.4546
$ws.Range("I20").Value = 1304.3462
$ws.Range("J20").Value = 2358.1667
$ws.Range("K20").Value = 1304.3462
$ws.Range("L20").Value = 2358.1667
$ws.Range("M20").Value = -1057.3462
$ws.Range("N20").Value = -2852.1667
$ws.Range("H22").Value = 349.16666
$ws.Range("J22").Value = 398
$ws.Range("L22").Value = 398
$ws.Range("N22").Value = -744
$ws.Range("H64").Value = 1312.8572
$ws.Range("I64").Value = 1198
$ws.Range("K64").Value = 1198
$ws.Range("M64").Value = -973
$ws.Range("H67").Value = 1312.8572
$ws.Range("I67").Value = 1198
$ws.Range("K67").Value = 1198
$ws.Range("M67").Value = -418
$ws.Range("H80").Value = 663.1539
$ws.Range("J80").Value = 943.6667
$ws.Range("L80").Value = 943.6667
$ws.Range("N80").Value = -2939.6667
$ws.Range("H83").Value = 663.1539
$ws.Range("J83").Value = 943.6667
$ws.Range("L83").Value = 4718.3335
$ws.Range("N83").Value = -14702.3335
$ws.Range("H94").Value = 263.5
$ws.Range("I94").Value = 199
$ws.Range("K94").Value = 199
$ws.Range("M94").Value = 252
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 10000
$ws.Range("K96").Value = 10000
$ws.Range("M96").Value = -7254
$ws.Range("H105").Value = 5651.5454
$ws.Range("I105").Value = 6036.4287
$ws.Range("J105").Value = 4978
$ws.Range("K105").Value = 6036.4287
$ws.Range("L105").Value = 4978
$ws.Range("M105").Value = -4289.4287
$ws.Range("N105").Value = -8472
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H134").Value = 3235.5688
$ws.Range("I134").Value = 2212
$ws.Range("K134").Value = 6636
$ws.Range("M134").Value = -4101

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4880.875
$ws.Range("I15").Value = 679.3333
$ws.Range("K15").Value = 679.3333
$ws.Range("M15").Value = -509.3333
$ws.Range("H31").Value = 2373.5625
$ws.Range("I31").Value = 1372.2646
$ws.Range("J31").Value = 3508.3667
$ws.Range("K31").Value = 1372.2646
$ws.Range("L31").Value = 3508.3667
$ws.Range("M31").Value = -1077.2646
$ws.Range("N31").Value = -4098.3667
$ws.Range("H34").Value = 2373.5625
$ws.Range("I34").Value = 1372.2646
$ws.Range("J34").Value = 3508.3667
$ws.Range("K34").Value = 1372.2646
$ws.Range("L34").Value = 3508.3667
$ws.Range("M34").Value = -1170.2646
$ws.Range("N34").Value = -3912.3667
$ws.Range("H99").Value = 4520.826
$ws.Range("I99").Value = 3450.2727
$ws.Range("J99").Value = 5502.1665
$ws.Range("K99").Value = 3450.2727
$ws.Range("L99").Value = 5502.1665
$ws.Range("M99").Value = -1952.2727
$ws.Range("N99").Value = -8498.166499999999
$ws.Range("H122").Value = 2297.8076
$ws.Range("I122").Value = 1605.8235
$ws.Range("K122").Value = 4817.470499999999
$ws.Range("M122").Value = -2367.470499999999
$ws.Range("H126").Value = 4520.826
$ws.Range("I126").Value = 3450.2727
$ws.Range("J126").Value = 5502.1665
$ws.Range("K126").Value = 10350.8181
$ws.Range("L126").Value = 16506.4995
$ws.Range("M126").Value = -7880.8181
$ws.Range("N126").Value = -21446.4995
$ws.Range("H132").Value = 4354.5864
$ws.Range("I132").Value = 3287.4119
$ws.Range("J132").Value = 5866.4165
$ws.Range("K132").Value = 9862.235700000001
$ws.Range("L132").Value = 17599.2495
$ws.Range("M132").Value = -7332.235700000001
$ws.Range("N132").Value = -22659.2495
$ws.Range("H134").Value = 4435.676
$ws.Range("I134").Value = 3627.65
$ws.Range("K134").Value = 10882.95
$ws.Range("M134").Value = -8347.950000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 909298.25
$ws.Range("I8").Value = 909298.25
$ws.Range("K8").Value = 2727894.75
$ws.Range("M8").Value = -2727755.75
$ws.Range("H11").Value = 4000622
$ws.Range("J11").Value = 1055.5
$ws.Range("L11").Value = 3166.5
$ws.Range("N11").Value = -3446.5
$ws.Range("H15").Value = 340
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 433.33334
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 1300.00002
$ws.Range("M15").Value = -460
$ws.Range("N15").Value = -1580.00002
$ws.Range("H60").Value = 1000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 3000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -3502
$ws.Range("H114").Value = 624.41174
$ws.Range("I114").Value = 472.2857
$ws.Range("J114").Value = 730.9
$ws.Range("K114").Value = 1416.8571
$ws.Range("L114").Value = 2192.7
$ws.Range("M114").Value = 1837.1429
$ws.Range("N114").Value = -8700.700000000001
$ws.Range("H121").Value = 2011.3636
$ws.Range("I121").Value = 1677.6666
$ws.Range("J121").Value = 2411.8
$ws.Range("K121").Value = 5032.9998
$ws.Range("L121").Value = 7235.400000000001
$ws.Range("M121").Value = -3722.9998
$ws.Range("N121").Value = -9855.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34486844
$ws.Range("J70").Value = 90913610
$ws.Range("L70").Value = 90913610
$ws.Range("N70").Value = -90914150
$ws.Range("H73").Value = 34486844
$ws.Range("J73").Value = 90913610
$ws.Range("L73").Value = 90913610
$ws.Range("N73").Value = -90915482
$ws.Range("H80").Value = 1180389.4
$ws.Range("I80").Value = 1004251.1
$ws.Range("J80").Value = 1432015.4
$ws.Range("K80").Value = 1004251.1
$ws.Range("L80").Value = 1432015.4
$ws.Range("M80").Value = -1003253.1
$ws.Range("N80").Value = -1434011.4
$ws.Range("H83").Value = 1180389.4
$ws.Range("I83").Value = 1004251.1
$ws.Range("J83").Value = 1432015.4
$ws.Range("K83").Value = 5021255.5
$ws.Range("L83").Value = 7160077
$ws.Range("M83").Value = -5016263.5
$ws.Range("N83").Value = -7170061
$ws.Range("H102").Value = 1269.75
$ws.Range("I102").Value = 1337.6875
$ws.Range("K102").Value = 1337.6875
$ws.Range("M102").Value = 284.3125
$ws.Range("H113").Value = 446770.66
$ws.Range("I113").Value = 681277.6
$ws.Range("J113").Value = 7070.125
$ws.Range("K113").Value = 681277.6
$ws.Range("L113").Value = 7070.125
$ws.Range("M113").Value = -679107.6
$ws.Range("N113").Value = -11410.125
$ws.Range("H122").Value = 9455.041999999999
$ws.Range("I122").Value = 8572.352999999999
$ws.Range("J122").Value = 11598.714
$ws.Range("K122").Value = 25717.059
$ws.Range("L122").Value = 34796.142
$ws.Range("M122").Value = -23267.059
$ws.Range("N122").Value = -39696.142
$ws.Range("H123").Value = 43197.8
$ws.Range("J123").Value = 43197.8
$ws.Range("L123").Value = 43197.8
$ws.Range("N123").Value = -48097.8
$ws.Range("H133").Value = 49998.4
$ws.Range("J133").Value = 49998.4
$ws.Range("L133").Value = 49998.4
$ws.Range("N133").Value = -60118.4
$ws.Range("H134").Value = 54999.668
$ws.Range("J134").Value = 54999.668
$ws.Range("L134").Value = 164999.004
$ws.Range("N134").Value = -170069.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5636.174
$ws.Range("I7").Value = 6176.2666
$ws.Range("J7").Value = 4623.5
$ws.Range("K7").Value = 6176.2666
$ws.Range("L7").Value = 4623.5
$ws.Range("M7").Value = -6064.2666
$ws.Range("N7").Value = -4847.5
$ws.Range("H40").Value = 3990
$ws.Range("I40").Value = 3538.0417
$ws.Range("J40").Value = 4976.091
$ws.Range("K40").Value = 3538.0417
$ws.Range("L40").Value = 4976.091
$ws.Range("M40").Value = -3402.0417
$ws.Range("N40").Value = -5248.091
$ws.Range("H55").Value = 1066.8182
$ws.Range("J55").Value = 2461.4285
$ws.Range("L55").Value = 2461.4285
$ws.Range("N55").Value = -2807.4285
$ws.Range("H61").Value = 4720.9165
$ws.Range("I61").Value = 2961.3333
$ws.Range("K61").Value = 2961.3333
$ws.Range("M61").Value = -2759.3333
$ws.Range("H82").Value = 3949.5557
$ws.Range("I82").Value = 3897.3333
$ws.Range("J82").Value = 3975.6667
$ws.Range("K82").Value = 3897.3333
$ws.Range("L82").Value = 3975.6667
$ws.Range("M82").Value = -3536.3333
$ws.Range("N82").Value = -4697.6667
$ws.Range("H85").Value = 3949.5557
$ws.Range("I85").Value = 3897.3333
$ws.Range("J85").Value = 3975.6667
$ws.Range("K85").Value = 3897.3333
$ws.Range("L85").Value = 3975.6667
$ws.Range("M85").Value = -2649.3333
$ws.Range("N85").Value = -6471.6667
$ws.Range("H113").Value = 4720.9165
$ws.Range("I113").Value = 2961.3333
$ws.Range("K113").Value = 2961.3333
$ws.Range("M113").Value = -791.3332999999998
$ws.Range("H122").Value = 721883.5600000001
$ws.Range("I122").Value = 535054.8
$ws.Range("K122").Value = 1605164.4
$ws.Range("M122").Value = -1602714.4
$ws.Range("H126").Value = 5636.174
$ws.Range("I126").Value = 6176.2666
$ws.Range("J126").Value = 4623.5
$ws.Range("K126").Value = 18528.7998
$ws.Range("L126").Value = 13870.5
$ws.Range("M126").Value = -16058.7998
$ws.Range("N126").Value = -18810.5
$ws.Range("H135").Value = 112383.5
$ws.Range("J135").Value = 112383.5
$ws.Range("L135").Value = 112383.5
$ws.Range("N135").Value = -122523.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 7500
$ws.Range("I23").Value = 5001
$ws.Range("J23").Value = 9999
$ws.Range("K23").Value = 5001
$ws.Range("L23").Value = 9999
$ws.Range("M23").Value = -4772
$ws.Range("N23").Value = -10457
$ws.Range("H81").Value = 8102.421
$ws.Range("I81").Value = 1004.0714
$ws.Range("K81").Value = 2008.1428
$ws.Range("M81").Value = -947.1428000000001
$ws.Range("H84").Value = 8102.421
$ws.Range("I84").Value = 1004.0714
$ws.Range("K84").Value = 10040.714
$ws.Range("M84").Value = -4736.714
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H125").Value = 52998.75
$ws.Range("J125").Value = 52998.75
$ws.Range("L125").Value = 52998.75
$ws.Range("N125").Value = -62838.75
